$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.9999999914890052,
    0.99529476263550665,
    0.97572954008513291,
    0.96715682158861949,
    0.95903430687608127,
    0.94157426306775893,
    0.93998085102697304,
    0.93757974530202848,
    0.93780492873471677,
    0.93898535704265718,
    0.93922200094414643,
    0.93992950782737705,
    0.94608310000620266,
    0.94263435225300229,
    0.94011855699050373,
    0.93737663994980269,
    0.93366897440136387,
    0.93256008371255428,
    0.99361000467946581,
    0.9864930645473875,
    0.98509459569416435,
    0.98383009255370246,
    0.96869417014934256,
    0.95567270074325683,
    0.94921561368978447,
    0.93897769746669468,
    0.93560643177887581,
    0.9235180129779601,
    0.91519724302901007,
    0.91115196619159888,
    0.91299042772060668,
    0.91489579463347492,
    0.91437579344280095
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
